$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Date Colours")

# --- Fix dates on existing rows 9 and 10 (Sheet1) ---
$ws1.Range("A9").Value = 44373
$ws1.Range("A10").Value = 44373

# --- Fix Link Label text on rows 12-14 (Granites Gold Mine, was colour-coded differently) ---
$ws1.Range("F12").Value = "Granites Gold Mine"
$ws1.Range("F13").Value = "Granites Gold Mine"
$ws1.Range("F14").Value = "Granites Gold Mine"

# --- Insert new row 15 (Q9 -> Q10 Close contact, Alpha) pushing the previous row 15 down to 16 ---
$ws1.Rows("15").Insert()

# Row 15 (new content)
$ws1.Range("A15").Value = 44374
$ws1.Range("B15").Value = "Q9"
$ws1.Range("C15").Value = "Q10"
$ws1.Range("D15").Value = "Queensland"
$ws1.Range("F15").Value = "Close contact"
$ws1.Range("G15").Value = "Alpha (B.1.1.7)"

# Row 16 (previously row 15, shifted down: T1 m -> T5 m50, Granites Gold Mine, Delta)
$ws1.Range("C16").Value = "T5 m50"
$ws1.Range("F16").Value = "Granites Gold Mine"

# --- Add new rows 17-22 ---
$ws1.Range("A17").Value = 44375
$ws1.Range("B17").Value = "T1 m"
$ws1.Range("C17").Value = "TN1 m30"
$ws1.Range("D17").Value = "Northern Territory"
$ws1.Range("F17").Value = "Granites Gold Mine"
$ws1.Range("G17").Value = "Delta (B.1.617.2)"

$ws1.Range("A18").Value = 44375
$ws1.Range("B18").Value = "T1 m"
$ws1.Range("C18").Value = "TQ1"
$ws1.Range("D18").Value = "Northern Territory"
$ws1.Range("F18").Value = "Granites Gold Mine"
$ws1.Range("G18").Value = "Delta (B.1.617.2)"

$ws1.Range("A19").Value = 44376
$ws1.Range("B19").Value = "T5 m50"
$ws1.Range("C19").Value = "T6"
$ws1.Range("D19").Value = "Northern Territory"
$ws1.Range("F19").Value = "Granites Gold Mine"
$ws1.Range("G19").Value = "Delta (B.1.617.2)"

$ws1.Range("A20").Value = 44376
$ws1.Range("B20").Value = "T5 m50"
$ws1.Range("C20").Value = "T7 w"
$ws1.Range("D20").Value = "Northern Territory"
$ws1.Range("F20").Value = "Granites Gold Mine"
$ws1.Range("G20").Value = "Delta (B.1.617.2)"

$ws1.Range("A21").Value = 44376
$ws1.Range("B21").Value = "T1 m"
$ws1.Range("C21").Value = "TQ2"
$ws1.Range("D21").Value = "Northern Territory"
$ws1.Range("F21").Value = "Granites Gold Mine"
$ws1.Range("G21").Value = "Delta (B.1.617.2)"

$ws1.Range("A22").Value = 44376
$ws1.Range("B22").Value = "? w19 "
$ws1.Range("C22").Value = "? w19 "
$ws1.Range("D22").Value = "Queensland"
$ws1.Range("F22").Value = "Prince Charles Hospital"
$ws1.Range("G22").Value = "Unknown"

# --- Apply date format (d-mmm) to column A for new rows so style matches ---
$ws1.Range("A15:A22").NumberFormat = "d-mmm"

# --- Resize Table1 to cover the new rows ---
$tbl1 = $ws1.ListObjects.Item("Table1")
$tbl1.Resize($ws1.Range("A1:G22"))

# --- Update Date Colours sheet: shift colour gradient down one row, inserting a new lighter shade at top ---
$ws2.Range("B2").Value = "#f4eaf4"
$ws2.Range("B3").Value = "#e9d6ea"
$ws2.Range("B4").Value = "#ddc1df"
$ws2.Range("B5").Value = "#d2add4"
$ws2.Range("B6").Value = "#c699ca"
$ws2.Range("B7").Value = "#ba85bf"
$ws2.Range("B8").Value = "#ae71b4"
$ws2.Range("B9").Value = "#a25daa"
$ws2.Range("B10").Value = "#96499f"
$ws2.Range("B11").Value = "#893395"

# --- Switch active sheet/selection from Date Colours back to Sheet1, cell A22 ---
$ws1.Activate()
$ws1.Range("A22").Select()
